$d = $word.ActiveDocument

# --- Step 1: the "_GoBack" bookmark currently sits at the very end of the
# document (right after "Implement some of the GUI"). It needs to move so
# that it sits in the middle of the "NEEDS SOME WORK" text (the location of
# the most recent edit). Remove it from its old spot first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: locate " NEEDS SOME WORK" and re-insert the "_GoBack" bookmark
# right after "NEEDS" (i.e. splitting " \u2013 NEEDS SOME WORK" into
# " \u2013 NEEDS" + bookmark + " SOME WORK").
$rngNeeds = $d.Content
$rngNeeds.Find.Execute("NEEDS SOME WORK")
$splitPos = $rngNeeds.Start + 5   # length of "NEEDS"
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Step 3: append " - DONE" as its own run to the "Implement some of the
# GUI" bullet (this item is now complete).
$rngGui = $d.Content
$rngGui.Find.Execute("Implement some of the GUI")
$endPos = $rngGui.End
$insertionPoint = $d.Range($endPos, $endPos)
$insertionPoint.InsertAfter(" - DONE")

# Force the appended text to live in its own <w:r> run (matching how Word
# keeps distinctly-authored text in separate runs) rather than being merged
# into the preceding run, by briefly bookmarking and then removing the
# bookmark at the split point.
$splitRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("zzTempRunSplit", $splitRange)
$d.Bookmarks("zzTempRunSplit").Delete()
